$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 253
$ws.Range("I2").Value = 749
$ws.Range("J2").Value = 2931
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 830
$ws.Range("M2").Value = 51
$ws.Range("N2").Value = 528
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 40
$ws.Range("S2").Value = 309
$ws.Range("T2").Value = 513
$ws.Range("U2").Value = 34
$ws.Range("V2").Value = 4575
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 4651
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 61
$ws.Range("AA2").Value = 21
